# Auto-generated edit script: apply scheduled market-data refresh
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N)
# across the ALC/ARM/BSM/CRP/CUL/LTW/WVR sheets per the commit's data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 2875  # H74: 2980 -> 2875
$ws.Cells.Item(74, 9).Value = 2800  # I74: 2950 -> 2800
$ws.Cells.Item(74, 11).Value = 2800  # K74: 2950 -> 2800
$ws.Cells.Item(74, 13).Value = -1864  # M74: -2014 -> -1864
$ws.Cells.Item(77, 8).Value = 2875  # H77: 2980 -> 2875
$ws.Cells.Item(77, 9).Value = 2800  # I77: 2950 -> 2800
$ws.Cells.Item(77, 11).Value = 14000  # K77: 14750 -> 14000
$ws.Cells.Item(77, 13).Value = -9320  # M77: -10070 -> -9320
$ws.Cells.Item(132, 8).Value = 10425490  # H132: 11373189 -> 10425490
$ws.Cells.Item(132, 9).Value = 13900342  # I132: 14717963 -> 13900342
$ws.Cells.Item(132, 10).Value = 934  # J132: 959.6 -> 934
$ws.Cells.Item(132, 11).Value = 41701026  # K132: 44153889 -> 41701026
$ws.Cells.Item(132, 12).Value = 2802  # L132: 2878.8 -> 2802
$ws.Cells.Item(132, 13).Value = -41698496  # M132: -44151359 -> -41698496
$ws.Cells.Item(132, 14).Value = -7862  # N132: -7938.8 -> -7862
$ws.Cells.Item(138, 8).Value = 1600.1613  # H138: 1767.6471 -> 1600.1613
$ws.Cells.Item(138, 10).Value = 3497.3333  # J138: 3497.8333 -> 3497.3333
$ws.Cells.Item(138, 12).Value = 10491.9999  # L138: 10493.4999 -> 10491.9999
$ws.Cells.Item(138, 14).Value = -20771.9999  # N138: -20773.4999 -> -20771.9999
$ws.Cells.Item(141, 8).Value = 1275.2222  # H141: 1246.9518 -> 1275.2222
$ws.Cells.Item(141, 9).Value = 1091.9066  # I141: 1066.1948 -> 1091.9066
$ws.Cells.Item(141, 11).Value = 3275.7198  # K141: 3198.5844 -> 3275.7198
$ws.Cells.Item(141, 13).Value = 1904.2802  # M141: 1981.4156 -> 1904.2802

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 20999.152  # H32: 1187.3 -> 20999.152
$ws.Cells.Item(32, 9).Value = 4546.7617  # I32: 1190.2291 -> 4546.7617
$ws.Cells.Item(32, 10).Value = 193749.25  # J32: 1117 -> 193749.25
$ws.Cells.Item(32, 11).Value = 4546.7617  # K32: 1190.2291 -> 4546.7617
$ws.Cells.Item(32, 12).Value = 193749.25  # L32: 1117 -> 193749.25
$ws.Cells.Item(32, 13).Value = -4259.7617  # M32: -903.2291 -> -4259.7617
$ws.Cells.Item(32, 14).Value = -194323.25  # N32: -1691 -> -194323.25
$ws.Cells.Item(45, 8).Value = 57630.832  # H45: 64704.188 -> 57630.832
$ws.Cells.Item(45, 9).Value = 68293  # I45: 78639 -> 68293
$ws.Cells.Item(45, 11).Value = 68293  # K45: 78639 -> 68293
$ws.Cells.Item(45, 13).Value = -67916  # M45: -78262 -> -67916
$ws.Cells.Item(132, 8).Value = 3793.2856  # H132: 3587.1226 -> 3793.2856
$ws.Cells.Item(132, 9).Value = 3700.5386  # I132: 3703.6924 -> 3700.5386
$ws.Cells.Item(132, 10).Value = 4999  # J132: 3132.5 -> 4999
$ws.Cells.Item(132, 11).Value = 11101.6158  # K132: 11111.0772 -> 11101.6158
$ws.Cells.Item(132, 12).Value = 14997  # L132: 9397.5 -> 14997
$ws.Cells.Item(132, 13).Value = -8571.6158  # M132: -8581.0772 -> -8571.6158
$ws.Cells.Item(132, 14).Value = -20057  # N132: -14457.5 -> -20057

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 82857.5  # H86: 72769.44 -> 82857.5
$ws.Cells.Item(86, 9).Value = 95975.5  # I86: 88747.08 -> 95975.5
$ws.Cells.Item(86, 10).Value = 4149.5  # J86: 3533 -> 4149.5
$ws.Cells.Item(86, 11).Value = 95975.5  # K86: 88747.08 -> 95975.5
$ws.Cells.Item(86, 12).Value = 4149.5  # L86: 3533 -> 4149.5
$ws.Cells.Item(86, 13).Value = -94852.5  # M86: -87624.08 -> -94852.5
$ws.Cells.Item(86, 14).Value = -6395.5  # N86: -5779 -> -6395.5
$ws.Cells.Item(89, 8).Value = 82857.5  # H89: 72769.44 -> 82857.5
$ws.Cells.Item(89, 9).Value = 95975.5  # I89: 88747.08 -> 95975.5
$ws.Cells.Item(89, 10).Value = 4149.5  # J89: 3533 -> 4149.5
$ws.Cells.Item(89, 11).Value = 479877.5  # K89: 443735.4 -> 479877.5
$ws.Cells.Item(89, 12).Value = 20747.5  # L89: 17665 -> 20747.5
$ws.Cells.Item(89, 13).Value = -474261.5  # M89: -438119.4 -> -474261.5
$ws.Cells.Item(89, 14).Value = -31979.5  # N89: -28897 -> -31979.5
$ws.Cells.Item(100, 8).Value = 15000  # H100: 25000 -> 15000
$ws.Cells.Item(100, 10).Value = 15000  # J100: 25000 -> 15000
$ws.Cells.Item(100, 12).Value = 15000  # L100: 25000 -> 15000
$ws.Cells.Item(100, 14).Value = -17164  # N100: -27164 -> -17164

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1297.9565  # H58: 1175.0317 -> 1297.9565
$ws.Cells.Item(58, 9).Value = 1108.35  # I58: 892.9039 -> 1108.35
$ws.Cells.Item(58, 10).Value = 2562  # J58: 2508.7273 -> 2562
$ws.Cells.Item(58, 11).Value = 1108.35  # K58: 892.9039 -> 1108.35
$ws.Cells.Item(58, 12).Value = 2562  # L58: 2508.7273 -> 2562
$ws.Cells.Item(58, 13).Value = -905.3499999999999  # M58: -689.9039 -> -905.3499999999999
$ws.Cells.Item(58, 14).Value = -2968  # N58: -2914.7273 -> -2968
$ws.Cells.Item(70, 8).Value = 11029.667  # H70: 16500 -> 11029.667
$ws.Cells.Item(70, 10).Value = 11029.667  # J70: 16500 -> 11029.667
$ws.Cells.Item(70, 12).Value = 11029.667  # L70: 16500 -> 11029.667
$ws.Cells.Item(70, 14).Value = -11659.667  # N70: -17130 -> -11659.667
$ws.Cells.Item(73, 8).Value = 11029.667  # H73: 16500 -> 11029.667
$ws.Cells.Item(73, 10).Value = 11029.667  # J73: 16500 -> 11029.667
$ws.Cells.Item(73, 12).Value = 11029.667  # L73: 16500 -> 11029.667
$ws.Cells.Item(73, 14).Value = -13213.667  # N73: -18684 -> -13213.667
$ws.Cells.Item(81, 8).Value = 36307.8  # H81: 37449.8 -> 36307.8
$ws.Cells.Item(81, 10).Value = 36307.8  # J81: 37449.8 -> 36307.8
$ws.Cells.Item(81, 12).Value = 36307.8  # L81: 37449.8 -> 36307.8
$ws.Cells.Item(81, 14).Value = -38303.8  # N81: -39445.8 -> -38303.8
$ws.Cells.Item(82, 8).Value = 29999  # H82: 0 -> 29999
$ws.Cells.Item(82, 10).Value = 29999  # J82: 0 -> 29999
$ws.Cells.Item(82, 12).Value = 29999  # L82: 0 -> 29999
$ws.Cells.Item(82, 14).Value = -30721  # N82: None -> -30721
$ws.Cells.Item(84, 8).Value = 36307.8  # H84: 37449.8 -> 36307.8
$ws.Cells.Item(84, 10).Value = 36307.8  # J84: 37449.8 -> 36307.8
$ws.Cells.Item(84, 12).Value = 108923.4  # L84: 112349.4 -> 108923.4
$ws.Cells.Item(84, 14).Value = -118907.4  # N84: -122333.4 -> -118907.4
$ws.Cells.Item(85, 8).Value = 29999  # H85: 0 -> 29999
$ws.Cells.Item(85, 10).Value = 29999  # J85: 0 -> 29999
$ws.Cells.Item(85, 12).Value = 29999  # L85: 0 -> 29999
$ws.Cells.Item(85, 14).Value = -32495  # N85: None -> -32495
$ws.Cells.Item(86, 8).Value = 3600.75  # H86: 3488.125 -> 3600.75
$ws.Cells.Item(86, 9).Value = 2750  # I86: 2666.6667 -> 2750
$ws.Cells.Item(86, 10).Value = 3884.3333  # J86: 3981 -> 3884.3333
$ws.Cells.Item(86, 11).Value = 2750  # K86: 2666.6667 -> 2750
$ws.Cells.Item(86, 12).Value = 3884.3333  # L86: 3981 -> 3884.3333
$ws.Cells.Item(86, 13).Value = -1627  # M86: -1543.6667 -> -1627
$ws.Cells.Item(86, 14).Value = -6130.3333  # N86: -6227 -> -6130.3333
$ws.Cells.Item(89, 8).Value = 3600.75  # H89: 3488.125 -> 3600.75
$ws.Cells.Item(89, 9).Value = 2750  # I89: 2666.6667 -> 2750
$ws.Cells.Item(89, 10).Value = 3884.3333  # J89: 3981 -> 3884.3333
$ws.Cells.Item(89, 11).Value = 13750  # K89: 13333.3335 -> 13750
$ws.Cells.Item(89, 12).Value = 19421.6665  # L89: 19905 -> 19421.6665
$ws.Cells.Item(89, 13).Value = -8134  # M89: -7717.333500000001 -> -8134
$ws.Cells.Item(89, 14).Value = -30653.6665  # N89: -31137 -> -30653.6665
$ws.Cells.Item(94, 8).Value = 999.5  # H94: 1032.5333 -> 999.5
$ws.Cells.Item(94, 9).Value = 751.1667  # I94: 800.8 -> 751.1667
$ws.Cells.Item(94, 10).Value = 1123.6666  # J94: 1148.4 -> 1123.6666
$ws.Cells.Item(94, 11).Value = 751.1667  # K94: 800.8 -> 751.1667
$ws.Cells.Item(94, 12).Value = 1123.6666  # L94: 1148.4 -> 1123.6666
$ws.Cells.Item(94, 13).Value = -300.1667  # M94: -349.8 -> -300.1667
$ws.Cells.Item(94, 14).Value = -2025.6666  # N94: -2050.4 -> -2025.6666
$ws.Cells.Item(136, 8).Value = 1297.9565  # H136: 1175.0317 -> 1297.9565
$ws.Cells.Item(136, 9).Value = 1108.35  # I136: 892.9039 -> 1108.35
$ws.Cells.Item(136, 10).Value = 2562  # J136: 2508.7273 -> 2562
$ws.Cells.Item(136, 11).Value = 3325.05  # K136: 2678.7117 -> 3325.05
$ws.Cells.Item(136, 12).Value = 7686  # L136: 7526.1819 -> 7686
$ws.Cells.Item(136, 13).Value = -775.0499999999997  # M136: -128.7116999999998 -> -775.0499999999997
$ws.Cells.Item(136, 14).Value = -12786  # N136: -12626.1819 -> -12786

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 162.25  # H2: 199.5 -> 162.25
$ws.Cells.Item(2, 9).Value = 74.5  # I2: 99 -> 74.5
$ws.Cells.Item(2, 10).Value = 250  # J2: 300 -> 250
$ws.Cells.Item(2, 11).Value = 447  # K2: 594 -> 447
$ws.Cells.Item(2, 12).Value = 1500  # L2: 1800 -> 1500
$ws.Cells.Item(2, 13).Value = -334  # M2: -481 -> -334
$ws.Cells.Item(2, 14).Value = -1726  # N2: -2026 -> -1726
$ws.Cells.Item(100, 8).Value = 2333.3333  # H100: 3200 -> 2333.3333
$ws.Cells.Item(100, 10).Value = 2333.3333  # J100: 3200 -> 2333.3333
$ws.Cells.Item(100, 12).Value = 6999.999899999999  # L100: 9600 -> 6999.999899999999
$ws.Cells.Item(100, 14).Value = -8621.999899999999  # N100: -11222 -> -8621.999899999999
$ws.Cells.Item(108, 8).Value = 2432.6667  # H108: 618 -> 2432.6667
$ws.Cells.Item(108, 9).Value = 2432.6667  # I108: 618 -> 2432.6667
$ws.Cells.Item(108, 11).Value = 7298.000100000001  # K108: 1854 -> 7298.000100000001
$ws.Cells.Item(108, 13).Value = -4418.000100000001  # M108: 1026 -> -4418.000100000001
$ws.Cells.Item(115, 8).Value = 3891.0908  # H115: 3435.0557 -> 3891.0908
$ws.Cells.Item(115, 9).Value = 882  # I115: 824.6667 -> 882
$ws.Cells.Item(115, 10).Value = 5019.5  # J115: 3957.1333 -> 5019.5
$ws.Cells.Item(115, 11).Value = 2646  # K115: 2474.0001 -> 2646
$ws.Cells.Item(115, 12).Value = 15058.5  # L115: 11871.3999 -> 15058.5
$ws.Cells.Item(115, 13).Value = -1471  # M115: -1299.0001 -> -1471
$ws.Cells.Item(115, 14).Value = -17408.5  # N115: -14221.3999 -> -17408.5
$ws.Cells.Item(120, 8).Value = 1000030  # H120: 7465 -> 1000030
$ws.Cells.Item(120, 9).Value = 1000030  # I120: 7465 -> 1000030
$ws.Cells.Item(120, 11).Value = 3000090  # K120: 22395 -> 3000090
$ws.Cells.Item(120, 13).Value = -2995252  # M120: -17557 -> -2995252

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 1012329.5  # H46: 506234.7 -> 1012329.5
$ws.Cells.Item(46, 9).Value = 377  # I46: 240.5 -> 377
$ws.Cells.Item(46, 10).Value = 2024282  # J46: 1265226 -> 2024282
$ws.Cells.Item(46, 11).Value = 377  # K46: 240.5 -> 377
$ws.Cells.Item(46, 12).Value = 2024282  # L46: 1265226 -> 2024282
$ws.Cells.Item(46, 13).Value = -189  # M46: -52.5 -> -189
$ws.Cells.Item(46, 14).Value = -2024658  # N46: -1265602 -> -2024658
$ws.Cells.Item(61, 8).Value = 1762.5883  # H61: 1350.0454 -> 1762.5883
$ws.Cells.Item(61, 9).Value = 1795  # I61: 1282.7778 -> 1795
$ws.Cells.Item(61, 10).Value = 1733.7778  # J61: 1396.6154 -> 1733.7778
$ws.Cells.Item(61, 11).Value = 1795  # K61: 1282.7778 -> 1795
$ws.Cells.Item(61, 12).Value = 1733.7778  # L61: 1396.6154 -> 1733.7778
$ws.Cells.Item(61, 13).Value = -1593  # M61: -1080.7778 -> -1593
$ws.Cells.Item(61, 14).Value = -2137.7778  # N61: -1800.6154 -> -2137.7778
$ws.Cells.Item(113, 8).Value = 1762.5883  # H113: 1350.0454 -> 1762.5883
$ws.Cells.Item(113, 9).Value = 1795  # I113: 1282.7778 -> 1795
$ws.Cells.Item(113, 10).Value = 1733.7778  # J113: 1396.6154 -> 1733.7778
$ws.Cells.Item(113, 11).Value = 1795  # K113: 1282.7778 -> 1795
$ws.Cells.Item(113, 12).Value = 1733.7778  # L113: 1396.6154 -> 1733.7778
$ws.Cells.Item(113, 13).Value = 375  # M113: 887.2221999999999 -> 375
$ws.Cells.Item(113, 14).Value = -6073.7778  # N113: -5736.6154 -> -6073.7778

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(70, 8).Value = 32000  # H70: 10500 -> 32000
$ws.Cells.Item(70, 10).Value = 32000  # J70: 10500 -> 32000
$ws.Cells.Item(70, 12).Value = 32000  # L70: 10500 -> 32000
$ws.Cells.Item(70, 14).Value = -32630  # N70: -11130 -> -32630
$ws.Cells.Item(73, 8).Value = 32000  # H73: 10500 -> 32000
$ws.Cells.Item(73, 10).Value = 32000  # J73: 10500 -> 32000
$ws.Cells.Item(73, 12).Value = 32000  # L73: 10500 -> 32000
$ws.Cells.Item(73, 14).Value = -34184  # N73: -12684 -> -34184
$ws.Cells.Item(132, 8).Value = 2642.3777  # H132: 1787.4286 -> 2642.3777
$ws.Cells.Item(132, 9).Value = 2680.5278  # I132: 1876.9246 -> 2680.5278
$ws.Cells.Item(132, 10).Value = 2489.7778  # J132: 1508.4117 -> 2489.7778
$ws.Cells.Item(132, 11).Value = 8041.5834  # K132: 5630.7738 -> 8041.5834
$ws.Cells.Item(132, 12).Value = 7469.3334  # L132: 4525.2351 -> 7469.3334
$ws.Cells.Item(132, 13).Value = -5511.5834  # M132: -3100.7738 -> -5511.5834
$ws.Cells.Item(132, 14).Value = -12529.3334  # N132: -9585.2351 -> -12529.3334
